# Scheduled runner update: refresh market-board price snapshots
# (currentAveragePrice / LevePrice / LeveProfit columns) across the
# per-class Anima Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 44500
$ws.Range("J3").Value = 44500
$ws.Range("L3").Value = 44500
$ws.Range("N3").Value = -44728
$ws.Range("H15").Value = 306.99
$ws.Range("I15").Value = 306.99
$ws.Range("K15").Value = 920.97
$ws.Range("M15").Value = -751.97
$ws.Range("H62").Value = 57733
$ws.Range("I62").Value = 50000
$ws.Range("J62").Value = 61599.5
$ws.Range("K62").Value = 50000
$ws.Range("L62").Value = 61599.5
$ws.Range("M62").Value = -49376
$ws.Range("N62").Value = -62847.5
$ws.Range("H65").Value = 57733
$ws.Range("I65").Value = 50000
$ws.Range("J65").Value = 61599.5
$ws.Range("K65").Value = 250000
$ws.Range("L65").Value = 307997.5
$ws.Range("M65").Value = -246880
$ws.Range("N65").Value = -314237.5
$ws.Range("H98").Value = 2370
$ws.Range("I98").Value = 2116.6667
$ws.Range("J98").Value = 2750
$ws.Range("K98").Value = 2116.6667
$ws.Range("L98").Value = 2750
$ws.Range("M98").Value = -618.6667000000002
$ws.Range("N98").Value = -5746
$ws.Range("H102").Value = 44500
$ws.Range("J102").Value = 44500
$ws.Range("L102").Value = 44500
$ws.Range("N102").Value = -50990
$ws.Range("H107").Value = 1861.6666
$ws.Range("I107").Value = 2800
$ws.Range("J107").Value = 923.3333
$ws.Range("K107").Value = 2800
$ws.Range("L107").Value = 923.3333
$ws.Range("M107").Value = -880
$ws.Range("N107").Value = -4763.3333
$ws.Range("H116").Value = 10650.714
$ws.Range("I116").Value = 14390
$ws.Range("K116").Value = 14390
$ws.Range("M116").Value = -10948
$ws.Range("H122").Value = 2370
$ws.Range("I122").Value = 2116.6667
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 6350.000100000001
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -3900.000100000001
$ws.Range("N122").Value = -13150
$ws.Range("H134").Value = 63260
$ws.Range("J134").Value = 63260
$ws.Range("L134").Value = 63260
$ws.Range("N134").Value = -73400
$ws.Range("H140").Value = 72174.92
$ws.Range("J140").Value = 72174.92
$ws.Range("L140").Value = 72174.92
$ws.Range("N140").Value = -82534.92

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10103939
$ws.Range("I61").Value = 14494786
$ws.Range("J61").Value = 4991.4
$ws.Range("K61").Value = 14494786
$ws.Range("L61").Value = 4991.4
$ws.Range("M61").Value = -14494574
$ws.Range("N61").Value = -5415.4
$ws.Range("H110").Value = 887.5
$ws.Range("I110").Value = 850
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 850
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1195
$ws.Range("N110").Value = -5090
$ws.Range("H136").Value = 10103939
$ws.Range("I136").Value = 14494786
$ws.Range("J136").Value = 4991.4
$ws.Range("K136").Value = 43484358
$ws.Range("L136").Value = 14974.2
$ws.Range("M136").Value = -43481808
$ws.Range("N136").Value = -20074.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 458.9091
$ws.Range("I94").Value = 408.66666
$ws.Range("J94").Value = 685
$ws.Range("K94").Value = 408.66666
$ws.Range("L94").Value = 685
$ws.Range("M94").Value = 42.33334000000002
$ws.Range("N94").Value = -1587
$ws.Range("H100").Value = 81071.25
$ws.Range("J100").Value = 81071.25
$ws.Range("L100").Value = 81071.25
$ws.Range("N100").Value = -83235.25
$ws.Range("H103").Value = 55989.332
$ws.Range("J103").Value = 55989.332
$ws.Range("L103").Value = 55989.332
$ws.Range("N103").Value = -58333.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5436.2056
$ws.Range("I31").Value = 1687.6086
$ws.Range("J31").Value = 7160.56
$ws.Range("K31").Value = 1687.6086
$ws.Range("L31").Value = 7160.56
$ws.Range("M31").Value = -1392.6086
$ws.Range("N31").Value = -7750.56
$ws.Range("H34").Value = 5436.2056
$ws.Range("I34").Value = 1687.6086
$ws.Range("J34").Value = 7160.56
$ws.Range("K34").Value = 1687.6086
$ws.Range("L34").Value = 7160.56
$ws.Range("M34").Value = -1485.6086
$ws.Range("N34").Value = -7564.56
$ws.Range("H74").Value = 20591.4
$ws.Range("J74").Value = 20591.4
$ws.Range("L74").Value = 20591.4
$ws.Range("N74").Value = -22339.4
$ws.Range("H77").Value = 20591.4
$ws.Range("J77").Value = 20591.4
$ws.Range("L77").Value = 61774.2
$ws.Range("N77").Value = -70510.20000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 858.0847
$ws.Range("I5").Value = 691.7632
$ws.Range("J5").Value = 1159.0476
$ws.Range("K5").Value = 2075.2896
$ws.Range("L5").Value = 3477.142800000001
$ws.Range("M5").Value = -1963.2896
$ws.Range("N5").Value = -3701.142800000001
$ws.Range("H63").Value = 1950
$ws.Range("I63").Value = 1950
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 5850
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -5101
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 1391134.6
$ws.Range("I64").Value = 1225
$ws.Range("J64").Value = 1669116.6
$ws.Range("K64").Value = 3675
$ws.Range("L64").Value = 5007349.800000001
$ws.Range("M64").Value = -3405
$ws.Range("N64").Value = -5007889.800000001
$ws.Range("H66").Value = 1950
$ws.Range("I66").Value = 1950
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 17550
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -13806
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 1391134.6
$ws.Range("I67").Value = 1225
$ws.Range("J67").Value = 1669116.6
$ws.Range("K67").Value = 3675
$ws.Range("L67").Value = 5007349.800000001
$ws.Range("M67").Value = -2739
$ws.Range("N67").Value = -5009221.800000001
$ws.Range("H87").Value = 3000
$ws.Range("H90").Value = 3000
$ws.Range("H103").Value = 1812.5
$ws.Range("I103").Value = 500
$ws.Range("J103").Value = 4000
$ws.Range("K103").Value = 1500
$ws.Range("L103").Value = 12000
$ws.Range("M103").Value = -621
$ws.Range("N103").Value = -13758
$ws.Range("H107").Value = 2209.5
$ws.Range("I107").Value = 407
$ws.Range("J107").Value = 5042
$ws.Range("K107").Value = 1221
$ws.Range("L107").Value = 15126
$ws.Range("M107").Value = 699
$ws.Range("N107").Value = -18966
$ws.Range("H131").Value = 4099.054
$ws.Range("J131").Value = 5424.2593
$ws.Range("L131").Value = 16272.7779
$ws.Range("N131").Value = -26352.7779
$ws.Range("H132").Value = 2529.762
$ws.Range("I132").Value = 2362
$ws.Range("J132").Value = 2732.842
$ws.Range("K132").Value = 21258
$ws.Range("L132").Value = 24595.578
$ws.Range("M132").Value = -18728
$ws.Range("N132").Value = -29655.578
$ws.Range("H135").Value = 858.0847
$ws.Range("I135").Value = 691.7632
$ws.Range("J135").Value = 1159.0476
$ws.Range("K135").Value = 6225.8688
$ws.Range("L135").Value = 10431.4284
$ws.Range("M135").Value = -3690.8688
$ws.Range("N135").Value = -15501.4284
$ws.Range("H137").Value = 34038.457
$ws.Range("I137").Value = 7534.6875
$ws.Range("K137").Value = 22604.0625
$ws.Range("M137").Value = -17504.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -971
$ws.Range("H98").Value = 59903.25
$ws.Range("J98").Value = 59903.25
$ws.Range("L98").Value = 59903.25
$ws.Range("N98").Value = -65893.25
$ws.Range("H122").Value = 1253.5
$ws.Range("I122").Value = 1140.1818
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 3420.5454
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -970.5454
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 41673496
$ws.Range("I132").Value = 62508308
$ws.Range("J132").Value = 3875.75
$ws.Range("K132").Value = 187524924
$ws.Range("L132").Value = 11627.25
$ws.Range("M132").Value = -187522394
$ws.Range("N132").Value = -16687.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 94500
$ws.Range("J62").Value = 94500
$ws.Range("L62").Value = 94500
$ws.Range("N62").Value = -95748
$ws.Range("H65").Value = 94500
$ws.Range("J65").Value = 94500
$ws.Range("L65").Value = 283500
$ws.Range("N65").Value = -289740
$ws.Range("H76").Value = 63929.332
$ws.Range("J76").Value = 84644
$ws.Range("L76").Value = 84644
$ws.Range("N76").Value = -85320
$ws.Range("H79").Value = 63929.332
$ws.Range("J79").Value = 84644
$ws.Range("L79").Value = 84644
$ws.Range("N79").Value = -86984
$ws.Range("H122").Value = 5188.609
$ws.Range("I122").Value = 1999.6666
$ws.Range("J122").Value = 5666.95
$ws.Range("K122").Value = 5998.9998
$ws.Range("L122").Value = 17000.85
$ws.Range("M122").Value = -3548.9998
$ws.Range("N122").Value = -21900.85

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 82500
$ws.Range("J68").Value = 82500
$ws.Range("L68").Value = 82500
$ws.Range("N68").Value = -84122
$ws.Range("H71").Value = 82500
$ws.Range("J71").Value = 82500
$ws.Range("L71").Value = 247500
$ws.Range("N71").Value = -255612
$ws.Range("H82").Value = 69140.5
$ws.Range("J82").Value = 69140.5
$ws.Range("L82").Value = 69140.5
$ws.Range("N82").Value = -69906.5
$ws.Range("H85").Value = 69140.5
$ws.Range("J85").Value = 69140.5
$ws.Range("L85").Value = 69140.5
$ws.Range("N85").Value = -71792.5
$ws.Range("H96").Value = 4800
$ws.Range("I96").Value = 4220
$ws.Range("J96").Value = 5525
$ws.Range("K96").Value = 4220
$ws.Range("L96").Value = 5525
$ws.Range("M96").Value = -2847
$ws.Range("N96").Value = -8271
$ws.Range("H97").Value = 98572
$ws.Range("J97").Value = 98572
$ws.Range("L97").Value = 98572
$ws.Range("N97").Value = -100554
